# CRM-1980 - Add upcountry, prepaid, postpaid, invoice & contract fields
# to the Partner summary details excel template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header (row 1) + placeholder (row 2) values --------------------
# Written column-by-column (header then placeholder) so the shared-string
# table is built up in the same interleaved order as the source template.
$ws.Range("Q1").Value  = "Upcountry"
$ws.Range("Q2").Value  = "{excel_data_line_item:upcountry}"

$ws.Range("R1").Value  = "Upcountry Rate"
$ws.Range("R2").Value  = "{excel_data_line_item:upcountry_rate}"

$ws.Range("S1").Value  = "Upcountry Max Distance Threshold"
$ws.Range("S2").Value  = "{excel_data_line_item:upcountry_max_distance_threshold}"

$ws.Range("T1").Value  = "Upcountry Approval"
$ws.Range("T2").Value  = "{excel_data_line_item:upcountry_approval}"

$ws.Range("U1").Value  = "Upcountry Approval Email"
$ws.Range("U2").Value  = "{excel_data_line_item:upcountry_approval_email}"

$ws.Range("V1").Value  = "Invoice Email To"
$ws.Range("V2").Value  = "{excel_data_line_item:invoice_email_to}"

$ws.Range("W1").Value  = "Invoice Email Cc"
$ws.Range("W2").Value  = "{excel_data_line_item:invoice_email_cc}"

$ws.Range("X1").Value  = "Invoice Email Bcc"
$ws.Range("X2").Value  = "{excel_data_line_item:invoice_email_bcc}"

$ws.Range("Y1").Value  = "PrePaid or Postpaid"
$ws.Range("Y2").Value  = "{excel_data_line_item:is_prepaid}"

$ws.Range("Z1").Value  = "PrePaid Amoun"
$ws.Range("Z2").Value  = "{excel_data_line_item:prepaid_amount_limit}"

$ws.Range("AA1").Value = "PrePaid Notification Amount"
$ws.Range("AA2").Value = "{excel_data_line_item:prepaid_notification_amount}"

$ws.Range("AB1").Value = "PostPaid Credit Period"
$ws.Range("AB2").Value = "{excel_data_line_item:postpaid_credit_period}"

$ws.Range("AC1").Value = "PostPaid Notification Limit"
$ws.Range("AC2").Value = "{excel_data_line_item:postpaid_notification_limit}"

$ws.Range("AD1").Value = "PostPaid Grace Period"
$ws.Range("AD2").Value = "{excel_data_line_item:postpaid_grace_period}"

# --- Re-apply the existing header / value formatting to the new cells --
# Row 1: the new header cells reuse the same bold style already used by
# the rest of row 1 (A1:P1). AC1/AD1 sit past the sheet's previous right
# edge, so their format has to be copied in explicitly.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AC1:AD1").PasteSpecial(-4122) | Out-Null

# Row 2: the placeholder cells reuse the style already used by O2:P2.
$ws.Range("O2").Copy() | Out-Null
$ws.Range("Q2:AD2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Column widths for the new / resized columns -----------------------
$ws.Columns.Item(16).ColumnWidth = 37.5
$ws.Columns.Item(17).ColumnWidth = 28
$ws.Columns.Item(18).ColumnWidth = 30.5
$ws.Columns.Item(19).ColumnWidth = 32

# --- Sheet view selection ------------------------------------------------
$ws.Range("H22").Select() | Out-Null

# --- Print / page setup --------------------------------------------------
$ws.PageSetup.Orientation = 1 | Out-Null
